# -----------------------------------------------------------------------
# Refresh the cryptocurrency price/volume snapshot on Sheet1 (columns B-E,
# rows 2-51). A handful of rows also swap places (B/C/D/E) because the
# underlying ranking reshuffled a little between runs.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ col letter = new value }, only columns that actually changed
$updates = @{
    2 = @{ "D" = '91.259.17'; "E" = '  +2.07%  ' }
    3 = @{ "D" = '3.160.00'; "E" = '  +2.79%  ' }
    4 = @{ "E" = '  +0.20%  ' }
    5 = @{ "D" = '238.89'; "E" = '  +1.67%  ' }
    6 = @{ "D" = '620.45'; "E" = '  +0.49%  ' }
    7 = @{ "D" = '1.13'; "E" = '  +7.15%  ' }
    8 = @{ "D" = '0.374'; "E" = '  +4.20%  ' }
    9 = @{ "E" = '  -0.08%  ' }
    10 = @{ "D" = '3.153.55'; "E" = '  +2.73%  ' }
    11 = @{ "D" = '0.745'; "E" = '  +5.93%  ' }
    12 = @{ "E" = '  +2.58%  ' }
    13 = @{ "D" = '0.0000247'; "E" = '  -0.22%  ' }
    14 = @{ "D" = '35.30'; "E" = '  +1.36%  ' }
    15 = @{ "D" = '5.59'; "E" = '  +4.68%  ' }
    16 = @{ "D" = '91.446.09'; "E" = '  +2.44%  ' }
    17 = @{ "D" = '3.743.06'; "E" = '  +2.91%  ' }
    18 = @{ "D" = '3.154.22'; "E" = '  +3.35%  ' }
    19 = @{ "D" = '3.74'; "E" = '  +0.45%  ' }
    20 = @{ "E" = '  +11.51%  ' }
    21 = @{ "D" = '6.02'; "E" = '  +12.10%  ' }
    22 = @{ "D" = '456.33'; "E" = '  +6.04%  ' }
    23 = @{ "E" = '  -3.70%  ' }
    24 = @{ "D" = '9.22'; "E" = '  +6.25%  ' }
    25 = @{ "D" = '6.02'; "E" = '  +8.75%  ' }
    26 = @{ "D" = '89.25'; "E" = '  +3.62%  ' }
    27 = @{ "D" = '12.07'; "E" = '  +3.78%  ' }
    28 = @{ "D" = '3.320.33'; "E" = '  +2.51%  ' }
    29 = @{ "D" = '0.999'; "E" = '  -0.05%  ' }
    30 = @{ "D" = '0.132'; "E" = '  +49.32%  ' }
    31 = @{ "B" = 'Stellar'; "C" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; "D" = '0.232'; "E" = '  +19.56%  ' }
    32 = @{ "B" = 'Cronos'; "C" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; "D" = '0.171'; "E" = '  +10.64%  ' }
    33 = @{ "D" = '9.38'; "E" = '  +4.22%  ' }
    34 = @{ "B" = 'Kaspa'; "C" = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; "D" = '0.170'; "E" = '  +13.80%  ' }
    35 = @{ "B" = 'Binance-PegBSC-USD'; "C" = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; "D" = '1.00'; "E" = '  -9.56%  ' }
    36 = @{ "B" = 'RenderToken'; "C" = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'; "D" = '7.69'; "E" = '  +9.13%  ' }
    37 = @{ "B" = 'EthereumClassic'; "C" = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; "D" = '26.53'; "E" = '  +4.42%  ' }
    38 = @{ "D" = '510.87'; "E" = '  +4.09%  ' }
    39 = @{ "E" = '  +4.44%  ' }
    40 = @{ "D" = '1.36'; "E" = '  +8.71%  ' }
    41 = @{ "B" = 'MantraDAO'; "C" = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'; "D" = '3.92'; "E" = '  +5.65%  ' }
    42 = @{ "B" = 'PolygonEcosystemToken'; "C" = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; "D" = '0.455'; "E" = '  +15.02%  ' }
    43 = @{ "D" = '3.45'; "E" = '  -4.23%  ' }
    44 = @{ "D" = '22.14'; "E" = '  +0.39%  ' }
    45 = @{ "E" = '  -0.07%  ' }
    46 = @{ "B" = 'ARBITRUM'; "C" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; "D" = '0.715'; "E" = '  +6.95%  ' }
    47 = @{ "B" = 'Stacks'; "C" = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; "D" = '1.95'; "E" = '  +5.56%  ' }
    48 = @{ "B" = 'Monero'; "C" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; "D" = '157.89'; "E" = '  +1.79%  ' }
    49 = @{ "D" = '1.38'; "E" = '  +6.43%  ' }
    50 = @{ "D" = '4.51'; "E" = '  +4.83%  ' }
    51 = @{ "D" = '44.06'; "E" = '  -0.66%  ' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $value = $updates[$row][$col]
        $addr = "$col$row"
        $range = $ws.Range($addr)

        # "Price" values such as 238.89 / 1.00 / 0.170 look numeric to Excel's
        # auto-detection and would otherwise be coerced to a Double (losing
        # trailing zeros / exact text). The source sheet stores them as plain
        # text, so force a text number format for single-dot numeric strings
        # in column D, write the value, then drop back to the default style so
        # no extra formatting is left behind on the cell.
        if ($col -eq "D" -and $value -match "^[+-]?\d+(\.\d+)?$") {
            $range.NumberFormat = "@"
            $range.Value = $value
            $range.Style = "Normal"
        } else {
            $range.Value = $value
        }
    }
}

Write-Output ("Updated {0} rows" -f $updates.Count)
